# Updated BGR model - 2025-08-15 23:45
#
# The "wind" worksheet lists wind-resource cost-class rows in pairs that sit
# next to each other (e.g. cost class 2 then cost class 3 for the same CF
# class). This edit re-orders several of those adjacent pairs so that the
# higher/alternate cost class comes first, swapping the process id (col C),
# the description (col D), the commodity id (col K) and the lcoe_class
# number (col P) between the two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wind")

function Swap-Rows {
    param([int]$RowA, [int]$RowB)

    foreach ($col in @("C", "D", "K", "P")) {
        $addrA = "$col$RowA"
        $addrB = "$col$RowB"
        $valA = $ws.Range($addrA).Value2
        $valB = $ws.Range($addrB).Value2
        $ws.Range($addrA).Value2 = $valB
        $ws.Range($addrB).Value2 = $valA
    }
}

# Pairs of rows whose content gets swapped.
Swap-Rows 4 5
Swap-Rows 13 14
Swap-Rows 15 16
Swap-Rows 19 20
Swap-Rows 27 28
Swap-Rows 47 48
